$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 data
$ws.Range("A14").Value = 45957
$ws.Range("B14").Value = 733
$ws.Range("C14").Value = 22
$ws.Range("D14").Value = 711

# Row 15 data
$ws.Range("A15").Value = 45954
$ws.Range("B15").Value = 595
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 576

# Update selection to reflect the new active cell / selection range
$ws.Range("A15:D15").Select()
